$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange

# Replace the whole paragraph with the first run's text, then append the
# remaining runs individually so each lands in its own <a:r> (matching the
# run-split in the target: "Der Versuch John " / "the" / " Ripper ... lassen").
$tr.Text = "Der Versuch John "
$null = $tr.InsertAfter("the")
$null = $tr.InsertAfter(" Ripper über eine Server und Clientartige Struktur auf mehreren Systemen gegen die selben Hashwerte strukturiert laufen zu lassen")

# New second paragraph, indented one level, per the diff.
$null = $tr.InsertAfter("`rWurde leider eingestellt")

$tr2 = $tr.Paragraphs(2, 1)
$tr2.IndentLevel = 2
